$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "305112005"
$ws.Range("B8").Value = "israel"
$ws.Range("C8").Value = "hadad"
$ws.Range("D8").Value = "1234"
$ws.Range("E8").Value = $true
